$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Actual Output" (E) and "Pass" (F) columns for the test cases
# that previously had no result recorded. Values are written in the same
# order the original author entered them (Jared's ticket-info test cases
# first, then the User login cases, then the Tech login cases) so that the
# shared-string table is built up in the same sequence.

$ws.Range("E17").Value = "Ticket information is shown"
$ws.Range("F17").Value = "Yes"

$ws.Range("E18").Value = "Comment was added to ticket"
$ws.Range("F18").Value = "Yes"

$ws.Range("E19").Value = 'No error is given and "Assigned Tech: undefined (undefined)" added to ticket info'
$ws.Range("F19").Value = "No"

$ws.Range("E20").Value = "Ticket status shows closed"
$ws.Range("F20").Value = "Yes"

$ws.Range("E21").Value = "Ticket status shows opened"
$ws.Range("F21").Value = "Yes"

$ws.Range("E22").Value = "Ticket now shows tech assgined"
$ws.Range("F22").Value = "Yes"

$ws.Range("E4").Value = "User logged in and directed to home page"
$ws.Range("F4").Value = "Yes"

$ws.Range("E5").Value = "User not logged in and error was given"
$ws.Range("F5").Value = "Yes"

$ws.Range("E2").Value = "Tech logged in and directed to home page"
$ws.Range("F2").Value = "Yes"

$ws.Range("E3").Value = "Tech not logged in and error was given"
$ws.Range("F3").Value = "Yes"

# Column widths were tweaked (column A narrowed, column E widened to fit the
# longer "Actual Output" text that was just added, etc). ColumnWidth values
# here are chosen so the saved internal width lands as close as possible to
# the authored widths.
$ws.Columns.Item(1).ColumnWidth = 10.451822916666666
$ws.Columns.Item(2).ColumnWidth = 37.877604166666664
$ws.Columns.Item(3).ColumnWidth = 37.451822916666664
$ws.Columns.Item(4).ColumnWidth = 57.736979166666664
$ws.Columns.Item(5).ColumnWidth = 74.59244791666667
$ws.Columns.Item(6).ColumnWidth = 3.8776041666666665

# Selection moved from E12 to E35 (off the used range, from scrolling/review).
$ws.Range("E35").Select() | Out-Null
